$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Save" header in H1, matching the formatting of the existing
# header cells (bold, centered, bordered style used by B1:G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the new "Save" data value in H2 (plain/default formatting, like the
# numeric value cells B2:G2 carry no explicit style).
$ws.Range("H2").Value = 1
